$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in newly-recorded "pending" / resolved values for the existing six
#    trading-instrument blocks (rows 3-54).
# ---------------------------------------------------------------------------

# Block 1 - 创业板50（159949）, rows 3-9
$ws.Range("C7").Value = "pending"
$ws.Range("C8").Value = "pending"
$ws.Range("E8").Value = "pending"
$ws.Range("C9").Value = "pending"

# Block 2 - 300ETF（510300）, rows 12-18
$ws.Range("C16").Value = "pending"
$ws.Range("E16").Value = "pending"
$ws.Range("C17").Value = "pending"
$ws.Range("E17").Value = "pending"
$ws.Range("C18").Value = "pending"

# Block 3 - 科创50（588000）, rows 21-27
$ws.Range("C25").Value = "pending"
$ws.Range("C26").Value = "pending"
$ws.Range("E26").Value = "pending"
$ws.Range("C27").Value = "pending"

# Block 4 - 证券ETF（512880）, rows 30-36
$ws.Range("C34").Value = 1
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = "1.155/1.173"
$ws.Range("E35").Value = 1
$ws.Range("C36").Value = "pending"

# Block 5 - 生物医药（512290）, rows 39-45
$ws.Range("D44").Value = "1.979/2.043"
$ws.Range("E44").Value = 1
$ws.Range("C45").Value = "pending"

# Block 6 - 银行ETF（512800）, rows 48-54
$ws.Range("C52").Value = "pending"
$ws.Range("C53").Value = "pending"
$ws.Range("C54").Value = "pending"

# ---------------------------------------------------------------------------
# 2. Append a brand-new seventh block (5GETF 515050) as rows 57-63, built by
#    copying the formatting of the existing block-4 template (rows 30-36)
#    and then overwriting the values.
# ---------------------------------------------------------------------------

$ws.Range("A30:R36").Copy() | Out-Null
$ws.Range("A57:R63").PasteSpecial(-4122) | Out-Null
$ws.Range("A30:R36").Copy() | Out-Null
$ws.Range("A57:R63").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B57").Value = "5GETF（515050）"
$ws.Range("B58").Value = 7800
$ws.Range("B59").Value = "翻倍"

$ws.Range("B61").Value = "1.195/1.124"
$ws.Range("C61").Value = 1

$ws.Range("B62").Value = "1.175/1.190"
$ws.Range("C62").Value = 1
$ws.Range("D62").Value = "1.247/1.251"
$ws.Range("E62").Value = 1

$ws.Range("B63").Value = "1.095/1.109"
$ws.Range("C63").Value = "pending"

# ---------------------------------------------------------------------------
# 3. Update the view: drop the scrolled-away top-left cell and move the
#    active selection to C11.
# ---------------------------------------------------------------------------

$ws.Range("C11").Select() | Out-Null
